$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 128791
$ws.Range("E2").Value = 1830
$ws.Range("F2").Value = 1830
$ws.Range("G2").Value = 1898
$ws.Range("H2").Value = 1473
$ws.Range("I2").Value = 1492
$ws.Range("J2").Value = -19
$ws.Range("K2").Value = 171222
$ws.Range("L2").Value = 115491
$ws.Range("M2").Value = 55732
$ws.Range("N2").Value = 55375
$ws.Range("O2").Value = 357
$ws.Range("P2").Value = 11550
$ws.Range("Q2").Value = -4575
$ws.Range("R2").Value = -5213
$ws.Range("S2").Value = 5031
$ws.Range("T2").Value = 3665
$ws.Range("U2").Value = -8240
$ws.Range("V2").Value = 38266
$ws.Range("W2").Value = 1.42
$ws.Range("X2").Value = 1.14
$ws.Range("Y2").Value = 2.62
$ws.Range("Z2").Value = 0.85
$ws.Range("AA2").Value = 207.23
$ws.Range("AB2").Value = 458.55
$ws.Range("AC2").Value = 495
$ws.Range("AD2").Value = 30.9
$ws.Range("AE2").Value = 20078
$ws.Range("AF2").Value = 0.76
$ws.Range("AG2").Value = 191
$ws.Range("AH2").Value = 1.25
$ws.Range("AI2").Value = 34.35
$ws.Range("AJ2").Value = 301641683

# Row 3
$ws.Range("D3").Value = 97144
$ws.Range("E3").Value = -15019
$ws.Range("F3").Value = -15019
$ws.Range("G3").Value = -14618
$ws.Range("H3").Value = -12121
$ws.Range("I3").Value = -12054
$ws.Range("J3").Value = -67
$ws.Range("K3").Value = 173016
$ws.Range("L3").Value = 130358
$ws.Range("M3").Value = 42657
$ws.Range("N3").Value = 42374
$ws.Range("O3").Value = 284
$ws.Range("P3").Value = 11550
$ws.Range("Q3").Value = 6203
$ws.Range("R3").Value = -10807
$ws.Range("S3").Value = 11037
$ws.Range("T3").Value = 6082
$ws.Range("U3").Value = 121
$ws.Range("V3").Value = 50333
$ws.Range("W3").Value = -15.46
$ws.Range("X3").Value = -12.48
$ws.Range("Y3").Value = -24.66
$ws.Range("Z3").Value = -7.04
$ws.Range("AA3").Value = 305.59
$ws.Range("AB3").Value = 347.81
$ws.Range("AC3").Value = -3995
$ws.Range("AD3").Value = -2.08
$ws.Range("AE3").Value = 15364
$ws.Range("AF3").Value = 0.54
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 301641683

# Row 4
$ws.Range("D4").Value = 104142
$ws.Range("E4").Value = -1472
$ws.Range("F4").Value = -1472
$ws.Range("G4").Value = -689
$ws.Range("H4").Value = -1388
$ws.Range("I4").Value = -1212
$ws.Range("J4").Value = -175
$ws.Range("K4").Value = 172175
$ws.Range("L4").Value = 109422
$ws.Range("M4").Value = 62753
$ws.Range("N4").Value = 62643
$ws.Range("O4").Value = 110
$ws.Range("P4").Value = 19506
$ws.Range("Q4").Value = -15548
$ws.Range("R4").Value = 1295
$ws.Range("S4").Value = 13917
$ws.Range("T4").Value = 2070
$ws.Range("U4").Value = -17618
$ws.Range("V4").Value = 53273
$ws.Range("W4").Value = -1.41
$ws.Range("X4").Value = -1.33
$ws.Range("Y4").Value = -2.31
$ws.Range("Z4").Value = -0.8
$ws.Range("AA4").Value = 174.37
$ws.Range("AB4").Value = 219.32
$ws.Range("AC4").Value = -382
$ws.Range("AD4").Value = -21.37
$ws.Range("AE4").Value = 15043
$ws.Range("AF4").Value = 0.54
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 442267277

# Row 5
$ws.Range("D5").Value = 79012
$ws.Range("E5").Value = -5242
$ws.Range("F5").Value = -5242
$ws.Range("G5").Value = -4639
$ws.Range("H5").Value = -3407
$ws.Range("I5").Value = -3388
$ws.Range("J5").Value = -20
$ws.Range("K5").Value = 138181
$ws.Range("L5").Value = 80207
$ws.Range("M5").Value = 57975
$ws.Range("N5").Value = 57896
$ws.Range("O5").Value = 79
$ws.Range("P5").Value = 19506
$ws.Range("Q5").Value = 5401
$ws.Range("R5").Value = -744
$ws.Range("S5").Value = -9567
$ws.Range("T5").Value = 1079
$ws.Range("U5").Value = 4321
$ws.Range("V5").Value = 42461
$ws.Range("W5").Value = -6.63
$ws.Range("X5").Value = -4.31
$ws.Range("Y5").Value = -5.62
$ws.Range("Z5").Value = -2.2
$ws.Range("AA5").Value = 138.35
$ws.Range("AB5").Value = 247.28
$ws.Range("AC5").Value = -766
$ws.Range("AD5").Value = -8.449999999999999
$ws.Range("AE5").Value = 13903
$ws.Range("AF5").Value = 0.47
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 442267277

# Row 6
$ws.Range("D6").Value = 52651
$ws.Range("E6").Value = -4093
$ws.Range("F6").Value = -4093
$ws.Range("G6").Value = -4976
$ws.Range("H6").Value = -3882
$ws.Range("I6").Value = -3879
$ws.Range("K6").Value = 142829
$ws.Range("L6").Value = 75366
$ws.Range("M6").Value = 67463
$ws.Range("N6").Value = 67384
$ws.Range("P6").Value = 31506
$ws.Range("Q6").Value = 1655
$ws.Range("R6").Value = 3693
$ws.Range("S6").Value = 504
$ws.Range("T6").Value = 690
$ws.Range("U6").Value = 965
$ws.Range("V6").Value = 29147
$ws.Range("W6").Value = -7.77
$ws.Range("X6").Value = -7.37
$ws.Range("Y6").Value = -6.19
$ws.Range("Z6").Value = -2.76
$ws.Range("AA6").Value = 111.72
$ws.Range("AB6").Value = 145.91
$ws.Range("AC6").Value = -693
$ws.Range("AD6").Value = -10.69
$ws.Range("AE6").Value = 11154
$ws.Range("AF6").Value = 0.66
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 630000000
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 71419
$ws.Range("E7").Value = -4288
$ws.Range("G7").Value = -10237
$ws.Range("H7").Value = -10158
$ws.Range("I7").Value = -10191
$ws.Range("K7").Value = 145690
$ws.Range("L7").Value = 87956
$ws.Range("M7").Value = 57734
$ws.Range("N7").Value = 57572
$ws.Range("P7").Value = 31509
$ws.Range("Q7").Value = -2299
$ws.Range("R7").Value = -1451
$ws.Range("S7").Value = 2907
$ws.Range("T7").Value = 1283
$ws.Range("U7").Value = -1875
$ws.Range("W7").Value = -6
$ws.Range("X7").Value = -14.22
$ws.Range("Y7").Value = -16.31
$ws.Range("Z7").Value = -7.04
$ws.Range("AA7").Value = 152.35
$ws.Range("AC7").Value = -1617
$ws.Range("AD7").Value = -4.11
$ws.Range("AE7").Value = 9529
$ws.Range("AF7").Value = 0.7
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 78332
$ws.Range("E8").Value = 928
$ws.Range("G8").Value = 607
$ws.Range("H8").Value = 476
$ws.Range("I8").Value = 440
$ws.Range("K8").Value = 148332
$ws.Range("L8").Value = 90155
$ws.Range("M8").Value = 58177
$ws.Range("N8").Value = 58015
$ws.Range("P8").Value = 31509
$ws.Range("Q8").Value = 2965
$ws.Range("R8").Value = -1749
$ws.Range("S8").Value = -1800
$ws.Range("T8").Value = 1343
$ws.Range("U8").Value = -310
$ws.Range("W8").Value = 1.18
$ws.Range("X8").Value = 0.61
$ws.Range("Y8").Value = 0.76
$ws.Range("Z8").Value = 0.32
$ws.Range("AA8").Value = 154.97
$ws.Range("AC8").Value = 70
$ws.Range("AD8").Value = 95.3
$ws.Range("AE8").Value = 9603
$ws.Range("AF8").Value = 0.6899999999999999
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 82964
$ws.Range("E9").Value = 1957
$ws.Range("G9").Value = 1295
$ws.Range("H9").Value = 1034
$ws.Range("I9").Value = 1007
$ws.Range("K9").Value = 149097
$ws.Range("L9").Value = 89949
$ws.Range("M9").Value = 59149
$ws.Range("N9").Value = 59027
$ws.Range("P9").Value = 31509
$ws.Range("Q9").Value = 2579
$ws.Range("R9").Value = -1688
$ws.Range("S9").Value = -1369
$ws.Range("T9").Value = 1369
$ws.Range("U9").Value = 810
$ws.Range("W9").Value = 2.36
$ws.Range("X9").Value = 1.25
$ws.Range("Y9").Value = 1.72
$ws.Range("Z9").Value = 0.7
$ws.Range("AA9").Value = 152.07
$ws.Range("AC9").Value = 160
$ws.Range("AD9").Value = 41.63
$ws.Range("AE9").Value = 9770
$ws.Range("AF9").Value = 0.68
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
